# Update LDLC prices history: insert a new "run" column (timestamp) right
# before the trailing "nom" / "url_produit" columns, carrying forward the
# latest known price (from the previous run column, AQ) for every product
# row that still has a price, and leaving it blank for rows that had no
# price in the last run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column AQ (43) held the most recent price snapshot and AR/AS (44/45) held
# the "nom"/"url_produit" columns. Insert a brand new column at AR (44) -
# this shifts the old AR -> AS and old AS -> AT, exactly like the diff.
$ws.Range("AR1").EntireColumn.Insert()

# New header for the freshly inserted timestamp column.
$ws.Range("AR1").Value = "2026-01-29 15:24:42"

# Figure out how many data rows exist (header is row 1).
$lastRow = $ws.UsedRange.Rows.Count

# Column AQ is now the previous run's price column (index 43); the new
# column we just inserted is AR (index 44). Copy the previous price
# forward for every row where it existed; leave the cell untouched
# (blank) otherwise, matching the source data exactly.
for ($r = 2; $r -le $lastRow; $r++) {
    $previousPrice = $ws.Cells.Item($r, 43).Value()
    if ($previousPrice -ne $null) {
        $ws.Cells.Item($r, 44).Value = $previousPrice
    }
}
